$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking Price values so they are not
# auto-converted to numbers by Excel (the source data keeps them as text).
$textCells = @("D5", "D6", "D11", "D12", "D16", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D30", "D31", "D36", "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D46", "D47")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin data (price, 1h volume change, and for three rows
# that got reordered/replaced: coin name + link as well).
$ws.Range("D2").Value = '65.246.22'
$ws.Range("E2").Value = '  -2.07%  '
$ws.Range("D3").Value = '3.477.71'
$ws.Range("E3").Value = '  -1.26%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '588.20'
$ws.Range("E5").Value = '  -2.92%  '
$ws.Range("D6").Value = '137.41'
$ws.Range("E6").Value = '  -4.30%  '
$ws.Range("D7").Value = '3.476.38'
$ws.Range("E7").Value = '  -1.29%  '
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("E9").Value = '  -3.10%  '
$ws.Range("E10").Value = '  -5.92%  '
$ws.Range("D11").Value = '7.17'
$ws.Range("E11").Value = '  -6.85%  '
$ws.Range("D12").Value = '0.382'
$ws.Range("E12").Value = '  -5.27%  '
$ws.Range("D13").Value = '4.071.58'
$ws.Range("E13").Value = '  -1.05%  '
$ws.Range("E14").Value = '  -6.78%  '
$ws.Range("D15").Value = '3.489.13'
$ws.Range("E15").Value = '  -1.58%  '
$ws.Range("D16").Value = '26.51'
$ws.Range("E16").Value = '  -7.51%  '
$ws.Range("E17").Value = '  -1.20%  '
$ws.Range("D18").Value = '65.142.02'
$ws.Range("E18").Value = '  -1.96%  '
$ws.Range("D19").Value = '9.71'
$ws.Range("E19").Value = '  -9.60%  '
$ws.Range("D20").Value = '5.77'
$ws.Range("E20").Value = '  -5.86%  '
$ws.Range("D21").Value = '13.87'
$ws.Range("E21").Value = '  -5.15%  '
$ws.Range("D22").Value = '388.96'
$ws.Range("E22").Value = '  -7.90%  '
$ws.Range("D23").Value = '0.554'
$ws.Range("E23").Value = '  -5.88%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = '72.60'
$ws.Range("E24").Value = '  -5.61%  '
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("B26").Value = 'WrappedeETH'
$ws.Range("C26").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D26").Value = '3.621.65'
$ws.Range("E26").Value = '  -1.31%  '
$ws.Range("B27").Value = 'LEO'
$ws.Range("C27").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D27").Value = '5.76'
$ws.Range("E27").Value = '  -0.19%  '
$ws.Range("E28").Value = '  -4.07%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("D30").Value = '7.29'
$ws.Range("E30").Value = '  -6.51%  '
$ws.Range("D31").Value = '8.18'
$ws.Range("E31").Value = '  -8.07%  '
$ws.Range("D33").Value = '3.500.30'
$ws.Range("E33").Value = '  -0.84%  '
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("E35").Value = '  -7.37%  '
$ws.Range("D36").Value = '23.10'
$ws.Range("E36").Value = '  -4.45%  '
$ws.Range("D37").Value = '171.85'
$ws.Range("E37").Value = '  -0.75%  '
$ws.Range("B38").Value = 'Fetch.AI'
$ws.Range("C38").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D38").Value = '1.19'
$ws.Range("E38").Value = '  -10.53%  '
$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").Value = '6.81'
$ws.Range("E39").Value = '  -9.55%  '
$ws.Range("E40").Value = '  -9.28%  '
$ws.Range("D41").Value = '4.72'
$ws.Range("E41").Value = '  -8.95%  '
$ws.Range("D42").Value = '0.0779'
$ws.Range("E42").Value = '  -3.25%  '
$ws.Range("D43").Value = '0.810'
$ws.Range("E43").Value = '  -4.78%  '
$ws.Range("D44").Value = '42.55'
$ws.Range("E44").Value = '  -6.67%  '
$ws.Range("D46").Value = '25.13'
$ws.Range("E46").Value = '  +9.85%  '
$ws.Range("D47").Value = '4.35'
$ws.Range("E47").Value = '  -12.44%  '
$ws.Range("E48").Value = '  +4.24%  '
$ws.Range("E49").Value = '  -8.11%  '
$ws.Range("E50").Value = '  -5.10%  '
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '2.221.77'
$ws.Range("E51").Value = '  -3.53%  '
